$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Step 1: In the paragraph "中雨，今天是农历五月初五，是中国传统节日：端午节，
# 这一天我们要吃粽子，赛龙舟。" merge the first two runs (which together read
# "中雨，今天是农历五月初五，是中国传统节日：端午节，这一天我们要吃粽子，赛龙舟")
# into a single run, while leaving the trailing "。" run untouched/separate.
# ---------------------------------------------------------------------------

$mergedText = "中雨，今天是农历五月初五，是中国传统节日：端午节，这一天我们要吃粽子，赛龙舟"

# Locate the paragraph that currently holds this text split across two runs.
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith($mergedText)) {
        $targetPara = $p
        break
    }
}

$fullSentence = $mergedText + "。"

# Touch the whole paragraph text via Find/Replace so the run-set gets
# re-normalised (this collapses the two runs that make up $mergedText, plus
# the trailing period run, into one single run).
$scoped = $targetPara.Range.Duplicate
$scoped.Find.Execute($fullSentence, $false, $false, $false, $false, $false, `
                      $true, 1, $false, $fullSentence, 2) | Out-Null

# Re-split the trailing "。" back into its own run by briefly toggling a
# character attribute on it and then reverting it. The toggle is enough to
# force a run boundary at that point without leaving any visible formatting
# difference behind.
$paraEnd = $targetPara.Range.End
$periodRange = $d.Range($paraEnd - 2, $paraEnd - 1)
$periodRange.Bold = 1
$periodRange2 = $d.Range($paraEnd - 2, $paraEnd - 1)
$periodRange2.Bold = 0

# ---------------------------------------------------------------------------
# Step 2: Insert a brand-new paragraph right after that paragraph containing
# the date "2022年6月7日星期二", matching the run layout Word already uses for
# the other date paragraphs in this document ("2" / "022" / rest).
# ---------------------------------------------------------------------------

# Find an existing date paragraph to copy the run/format layout from.
$dateTemplate = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("2022") -and $p.Range.Text.Contains("年6月3日星期五")) {
        $dateTemplate = $p
        break
    }
}

$srcRange = $dateTemplate.Range.Duplicate
$insertAt = $targetPara.Range.End
$insertPoint = $d.Range($insertAt, $insertAt)
$insertPoint.FormattedText = $srcRange.FormattedText

# Locate the paragraph we just inserted (immediately after $targetPara) and
# fix up its date text.
$newDatePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -eq $insertAt) {
        $newDatePara = $p
        break
    }
}
$dateScoped = $newDatePara.Range.Duplicate
$dateScoped.Find.Execute("年6月3日星期五", $false, $false, $false, $false, $false, `
                          $true, 1, $false, "年6月7日星期二", 2) | Out-Null

# ---------------------------------------------------------------------------
# Step 3: The (previously empty) paragraph that follows now gets a new run of
# text: "晴，今天是高考第一天，上午考语文，下午考数学。"
# ---------------------------------------------------------------------------

$emptyPara1 = $newDatePara.Next()
$emptyPara1.Range.Text = "晴，今天是高考第一天，上午考语文，下午考数学。"

# ---------------------------------------------------------------------------
# Step 4: The next (also previously empty) paragraph becomes a totally bare
# paragraph (no paragraph properties at all).
# ---------------------------------------------------------------------------

$emptyPara2 = $emptyPara1.Next()
$emptyPara2.Range.ParagraphFormat.Reset()
$emptyPara2.Format.Reset()

# ---------------------------------------------------------------------------
# Step 5: Append one more brand-new, totally bare empty paragraph at the very
# end of the document body.
# ---------------------------------------------------------------------------

$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()
